$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tier 1_obs")

$values = @{
    6  = "IpiNivo_Complete_PBMC"
    7  = "IpiNivo_Complete_kidney"
    8  = "IpiNivo_Complete_tumor_Center"
    9  = "IpiNivo_Complete_tumor_Far"
    10 = "IpiNivo_Complete_tumor_Near"
    11 = "IpiNivo_Mixed_PBMC"
    12 = "IpiNivo_Mixed_kidney"
    13 = "IpiNivo_Mixed_tumor_Center"
    14 = "IpiNivo_Mixed_tumor_Far"
    15 = "IpiNivo_Mixed_tumor_Near"
    16 = "IpiNivo_Resistant_PBMC"
    17 = "IpiNivo_Resistant_kidney"
    18 = "IpiNivo_Resistant_lymph_node"
    19 = "IpiNivo_Resistant_tumor_Center"
    20 = "IpiNivo_Resistant_tumor_Far"
    21 = "IpiNivo_Resistant_tumor_Near"
    22 = "NivoExposed_PBMC"
    23 = "NivoExposed_kidney"
    24 = "NivoExposed_tumor_Center"
    25 = "NivoExposed_tumor_Far"
    26 = "NivoExposed_tumor_Near"
    27 = "UT1_kidney"
    28 = "UT1_tumor_Center"
    29 = "UT1_tumor_Far"
    30 = "UT1_tumor_Near"
    31 = "UT2_kidney"
    32 = "UT2_tumor_Center"
    33 = "UT2_tumor_Far"
    34 = "UT2_tumor_Near"
}

foreach ($row in $values.Keys) {
    $ws.Range("B$row").Value = $values[$row]
}
